$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '22.462.58'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +0.01%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.571.51'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +0.09%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  -0.03%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '1.001'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.02%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '288.29'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.62%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.3721'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +0.93%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '48.27'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -4.03%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.3318'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -1.97%  '

$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.07489'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -0.37%  '

$ws.Range('B11').Value = 'Polygon'
$ws.Range('C11').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '1.132'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -1.18%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.001'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +0.07%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '20.70'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -2.17%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.941'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -1.19%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '6.894'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -1.40%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '1.568.09'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -0.15%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.00001117'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -0.14%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '87.77'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -2.60%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06740'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -0.24%  '

$ws.Range('E20').Value = '  +0.10%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.351'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -0.16%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '16.50'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +1.15%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '12.08'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -0.30%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '22.460.02'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +0.05%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.389'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +0.25%  '

$ws.Range('E26').Value = '  -3.07%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '153.03'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +2.28%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '19.69'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -1.10%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '5.018'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -0.70%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '124.24'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -0.38%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.744.72'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -0.03%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.051'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -0.57%  '

$ws.Range('E33').Value = '  -0.35%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '6.132'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -1.63%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '9.772'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -0.25%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.08313'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -0.58%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.02465'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -0.33%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.2272'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -0.90%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.06397'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -1.80%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '5.369'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -0.71%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.291'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -3.98%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '11.31'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +0.67%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.6303'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +1.42%  '

$ws.Range('E44').Value = '  +0.03%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '13.79'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -1.82%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.6157'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +5.23%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '3.778'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -0.08%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.054'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -0.32%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '125.41'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.22%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.211'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -1.96%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.07221'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -1.08%  '

